# Hortaliza, Vega Modelo de Temuco - Camote
# Insert a new weekly price record above row 43, pushing every
# subsequent record down by one row (row 115 -> row 116).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 43; everything from the old
# row 43 down (through 115) shifts down to 44..116, carrying its
# formatting (incl. the date number format on column D) with it.
$ws.Rows.Item(43).Insert()

# The row directly below (44) now holds what used to be row 43, so
# we can copy the constant, per-market/category columns from there
# into the freshly inserted row.
$ws.Range("A43").Value2 = $ws.Range("A44").Value2
$ws.Range("B43").Value2 = $ws.Range("B44").Value2
$ws.Range("C43").Value2 = $ws.Range("C44").Value2
$ws.Range("E43").Value2 = $ws.Range("E44").Value2
$ws.Range("F43").Value2 = $ws.Range("F44").Value2
$ws.Range("G43").Value2 = $ws.Range("G44").Value2
$ws.Range("H43").Value2 = $ws.Range("H44").Value2
$ws.Range("I43").Value2 = $ws.Range("I44").Value2
$ws.Range("J43").Value2 = $ws.Range("J44").Value2
$ws.Range("K43").Value2 = $ws.Range("K44").Value2
$ws.Range("L43").Value2 = $ws.Range("L44").Value2
$ws.Range("M43").Value2 = $ws.Range("M44").Value2
$ws.Range("R43").Value2 = $ws.Range("R44").Value2

# New record's own data.
$ws.Range("D43").Value2 = 44797
$ws.Range("N43").Value2 = "$/malla 20 kilos"
$ws.Range("O43").Value2 = "Perú"
$ws.Range("P43").Value2 = 1000
$ws.Range("Q43").Value2 = 20
